$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1.75
$ws.Range("F7").Value = 2.5

$ws.Range("T3").Value = 2.25

$ws.Range("Q4").Value = 44551
$ws.Range("Q4").NumberFormat = $ws.Range("Q3").NumberFormat
$ws.Range("R4").Value = "8.45 - 12.00"
$ws.Range("T4").Value = 3.25
$ws.Range("U4").Value = "Worked on character controller"

$ws.Range("F10").Formula = "=SUM(F3:F7)"

$ws.Range("U5").Select()
